$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing data (D:K) to (E:L)
$ws.Columns("D:D").Insert()

# Copy number formats/styles from column E (the old D, already shifted) into the
# new column D, one contiguous block at a time (matches the three financial
# statement tables: Income Statement, Balance Sheet, Cash Flow Statement)
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Populate new column D (latest FY period: 2018-12-31) with data
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 79000
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 20900
$ws.Range("D18").Value = 58100
$ws.Range("D20").Value = -37000
$ws.Range("D21").Value = 23500
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 21100
$ws.Range("D24").Value = 4700
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 16400
$ws.Range("D27").Value = 16400
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 37000
$ws.Range("D33").Value = 16400
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 16400
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 20600
$ws.Range("D42").Value = 27300
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 24600
$ws.Range("D49").Value = 16200
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 9900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 2118800
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1881800
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 237000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 16400
$ws.Range("D83").Value = 2400
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 24600
$ws.Range("D91").Value = -3300
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -45700
$ws.Range("D96").Value = -4600
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 20800
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -300
